# Pasajeros_AICM.xlsx update:
#  - Insert a new "Ago." (August) 2025 data row at the top of the monthly
#    table (row 6), pushing every existing row down by one.
#  - Resize the Excel Table ("Tabla3") / its AutoFilter so it keeps covering
#    exactly the data rows (B5:D97 instead of B5:D96).
#  - Bump the "Actualización: ..." footer note from Julio to Agosto 2025.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tabla3")

# Physically insert a blank row above the current first data row (row 6,
# "2025 / Jul."). This shifts all the rows below (including the footer rows)
# down by one and keeps the sheet dimension in sync automatically.
$ws.Rows.Item(6).Insert()

# Clone the formatting of the (now shifted) "2025 / Jun." row just below so
# the new row picks up the correct alternating row style used by the table.
$ws.Range("B8:D8").Copy()
$ws.Range("B6:D6").PasteSpecial(-4122)   # xlPasteFormats

# Write the new August 2025 figures into the freshly inserted row.
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "Ago."
$ws.Range("D6").Value = 3936.602

# The table ("Tabla3") used to span B5:D96; grow it by one row so the newly
# inserted row (and the row that slid down from the bottom of the old range)
# stay inside the table/autofilter range -> B5:D97.
$tbl.Resize($ws.Range("B5:D97"))

# Update the "Actualización: ..." note (now on row 98 after the insert).
$ws.Range("B98").Value = "Actualización: Agosto 2025."
